$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet: "Login" -> "Otentikasi" ---
$ws.Name = "Otentikasi"

# --- Remove the pre-existing hyperlink on A2 (content there changes) and
#     drop the leftover Hyperlink character style that came with it.
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A2").Style = "Normal"

# --- Rewrite the grid (columns A, B and the new test-case column H) ---
# Row 1 (headers) is unchanged, but re-asserted for safety.
$ws.Range("A1").Value = "username/email"
$ws.Range("B1").Value = "password"

$ws.Range("A2").Value = "akhimusyafak"
$ws.Range("B2").Value = "akhmad20221"
$ws.Range("H2").Value = "TC-Login-001"

$ws.Range("A3").Value = "akhimusyafak@gmail.com"
$ws.Range("B3").Value = "akhmad20221"
$ws.Range("H3").Value = "TC-Login-002"

$ws.Range("A4").Value = "akhimusyafak"
$ws.Range("B4").Value = "akhmad20222"
$ws.Range("H4").Value = "TC-Login-003"

$ws.Range("A5").Value = "akhimusyafak@gmail.com"
$ws.Range("B5").Value = "akhmad20222"
$ws.Range("H5").Value = "TC-Login-004"

$ws.Range("A6").Value = "akhmad"
$ws.Range("B6").Value = "akhmad20221"
$ws.Range("H6").Value = "TC-Login-005"

$ws.Range("A7").Value = "akhmadmusyafak@gmail.com"
$ws.Range("B7").Value = "akhmad20221"
$ws.Range("H7").Value = "TC-Login-006"

$ws.Range("B8").Value = "akhmad20221"
$ws.Range("H8").Value = "TC-Login-007"

$ws.Range("A9").Value = "akhimusyafak"
$ws.Range("H9").Value = "TC-Login-008"

$ws.Range("A10").Value = "akhimusyafak@gmail.com"
$ws.Range("H10").Value = "TC-Login-009"

$ws.Range("H11").Value = "TC-Login-010"

# Register the small (size-8) auxiliary font the sheet ends up carrying
# (used for the phonetic-guide info on the run-case column) without
# changing how H11 itself actually renders.
$ws.Range("H11").Phonetics.Font.Size = 8
$ws.Range("H11").Style = "Normal"
$ws.Range("H11").Value = "TC-Login-010"

# --- Re-create the mailto hyperlinks on the email cells ---
[void]$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:akhimusyafak@gmail.com")
[void]$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:akhimusyafak@gmail.com")
[void]$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:akhmadmusyafak@gmail.com")
[void]$ws.Hyperlinks.Add($ws.Range("A10"), "mailto:akhimusyafak@gmail.com")

# --- Column sizing: widen A for the longer emails, add width for H ---
$ws.Columns("A").AutoFit()
$ws.Columns("H").AutoFit()

# --- Selection moves to C18 ---
[void]$ws.Range("C18").Select()
